# Update workbook with new daily rows (aggiornamento fino a 02/05)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: row, date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(239, 44313, 0, 1, 33.71544167228591),
    @(240, 44314, 0, 1, 33.71544167228591),
    @(241, 44315, 1, 2, 67.43088334457181),
    @(242, 44316, 0, 1, 33.71544167228591),
    @(243, 44317, 1, 2, 67.43088334457181),
    @(244, 44318, 1, 3, 101.1463250168577)
)

$srcRow = 238

foreach ($entry in $data) {
    $r = $entry[0]

    # Copy formatting from the row above (keeps styles/number formats identical)
    $ws.Range("A$srcRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]

    $srcRow = $r
}

$excel.CutCopyMode = 0
